$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.558.86"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "2.624.19"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.26"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.42"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.54"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.374"
$ws.Range("E11").Value = "  +1.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.155"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("D13").Value = "3.090.43"
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.18"
$ws.Range("E14").Value = "  +12.09%  "
$ws.Range("D15").Value = "60.535.66"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").Value = "2.633.80"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.50"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E21").Value = "  -1.49%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.532"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.87"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("E27").Value = "  +5.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.02"
$ws.Range("E28").Value = "  +10.90%  "
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("E30").Value = "  +5.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.48"
$ws.Range("E31").Value = "  +4.70%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  +9.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.41"
$ws.Range("E35").Value = "  +4.41%  "
$ws.Range("E36").Value = "  +7.90%  "
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "331.33"
$ws.Range("E38").Value = "  +12.65%  "
$ws.Range("E39").Value = "  +4.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.28"
$ws.Range("E40").Value = "  +1.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.873"
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("E42").Value = "  +6.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.72"
$ws.Range("E43").Value = "  +4.08%  "
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "132.77"
$ws.Range("E45").Value = "  -3.42%  "
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0557"
$ws.Range("E48").Value = "  +2.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.609"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.75"
$ws.Range("E51").Value = "  +0.39%  "
